{"js": "// Remove the entire list-item paragraph \"Wait for a decision.\" (it\n// directly follows the \"Go to the court hearing, if necessary.\" item).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst target = paragraphs.items.find(\n  (p) => p.text.trim() === \"Wait for a decision.\"\n);\n\nif (target) {\n  target.delete();\n  await context.sync();\n}\n", "ps1": "# Remove the entire list-item paragraph \"Wait for a decision.\" (it\n# directly follows the \"Go to the court hearing, if necessary.\" item).\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Trim() -eq \"Wait for a decision.\") {\n        $p.Range.Delete()\n        break\n    }\n}\n"}
